# Update: Threat Alert Report - 2026-01-18 01:00
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'26-JAN-26"
$ws.Range("B2").Value = "SM-438"
$ws.Range("C2").Value = "Nile Air NP-232"
$ws.Range("D2").Value = 354
$ws.Range("E2").Value = 583
$ws.Range("F2").Value = -229
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "LOW THREAT"
$ws.Range("K2").Value = "SAR"

$ws.Range("A3").Value = "'26-JAN-26"
$ws.Range("B3").Value = "SM-438"
$ws.Range("C3").Value = "Nesma Airlines NE-153"
$ws.Range("D3").Value = 371
$ws.Range("E3").Value = 583
$ws.Range("F3").Value = -212
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("K3").Value = "SAR"

$ws.Range("A4").Value = "'26-JAN-26"
$ws.Range("B4").Value = "SM-438"
$ws.Range("C4").Value = "flynas XY-894"
$ws.Range("D4").Value = 509
$ws.Range("E4").Value = 583
$ws.Range("F4").Value = -74
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = -10
$ws.Range("J4").Value = "LOW THREAT"
$ws.Range("K4").Value = "SAR"

$ws.Range("A5").Value = "'26-JAN-26"
$ws.Range("B5").Value = "SM-438"
$ws.Range("C5").Value = "flynas XY-854"
$ws.Range("D5").Value = 509
$ws.Range("E5").Value = 583
$ws.Range("F5").Value = -74
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = -10
$ws.Range("J5").Value = "LOW THREAT"
$ws.Range("K5").Value = "SAR"

$ws.Range("A6").Value = "'30-JAN-26"
$ws.Range("B6").Value = "SM-438"
$ws.Range("C6").Value = "Nesma Airlines NE-151"
$ws.Range("D6").Value = 371
$ws.Range("E6").Value = 583
$ws.Range("F6").Value = -212
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "LOW THREAT"
$ws.Range("K6").Value = "SAR"

$ws.Range("A7").Value = "'16-FEB-26"
$ws.Range("B7").Value = "SM-438"
$ws.Range("C7").Value = "flyadeal F3-911"
$ws.Range("D7").Value = 369
$ws.Range("E7").Value = 826
$ws.Range("F7").Value = -457
$ws.Range("G7").Value = 15
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 15
$ws.Range("J7").Value = "LOW THREAT"
$ws.Range("K7").Value = "SAR"

$ws.Range("A8").Value = "'16-FEB-26"
$ws.Range("B8").Value = "SM-438"
$ws.Range("C8").Value = "Nile Air NP-232"
$ws.Range("D8").Value = 563
$ws.Range("E8").Value = 826
$ws.Range("F8").Value = -263
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "LOW THREAT"
$ws.Range("K8").Value = "SAR"

$ws.Range("A9").Value = "'16-FEB-26"
$ws.Range("B9").Value = "SM-438"
$ws.Range("C9").Value = "flynas XY-894"
$ws.Range("D9").Value = 629
$ws.Range("E9").Value = 826
$ws.Range("F9").Value = -197
$ws.Range("G9").Value = 40
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = -10
$ws.Range("J9").Value = "LOW THREAT"
$ws.Range("K9").Value = "SAR"

$ws.Range("A10").Value = "'16-FEB-26"
$ws.Range("B10").Value = "SM-438"
$ws.Range("C10").Value = "flynas XY-854"
$ws.Range("D10").Value = 629
$ws.Range("E10").Value = 826
$ws.Range("F10").Value = -197
$ws.Range("G10").Value = 40
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = -10
$ws.Range("J10").Value = "LOW THREAT"
$ws.Range("K10").Value = "SAR"

$ws.Range("A11").Value = "'20-FEB-26"
$ws.Range("B11").Value = "SM-438"
$ws.Range("C11").Value = "flynas XY-894"
$ws.Range("D11").Value = 599
$ws.Range("E11").Value = 686
$ws.Range("F11").Value = -87
$ws.Range("G11").Value = 40
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = -10
$ws.Range("J11").Value = "LOW THREAT"
$ws.Range("K11").Value = "SAR"

$ws.Range("A12").Value = "'06-MAR-26"
$ws.Range("B12").Value = "SM-438"
$ws.Range("C12").Value = "flynas XY-894"
$ws.Range("D12").Value = 949
$ws.Range("E12").Value = 1026
$ws.Range("F12").Value = -77
$ws.Range("G12").Value = 40
$ws.Range("H12").Value = 30
$ws.Range("I12").Value = -10
$ws.Range("J12").Value = "LOW THREAT"
$ws.Range("K12").Value = "SAR"

# Remove the now-obsolete trailing rows (old rows 13-15)
$ws.Rows("13:15").Delete()
